# cambios en botones y estructura
# Update the "Chilefilms" (Matriz/Individual) column C values on the
# "Estado" (Estado de Situacion Financiera) sheet to match the refreshed
# financial data. A few cells are cleared back to blank (no value).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Estado")

$ws.Range("C6").Value = 0
$ws.Range("C7").Value = 3255872574.56
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 2982274761.46
$ws.Range("C10").Value = 60073559317.1884
$ws.Range("C11").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("C15").Value = $null
$ws.Range("C16").Value = 10397155079.081
$ws.Range("C20").Value = 5414956957.90246
$ws.Range("C21").Value = 1080749265
$ws.Range("C22").Value = 2322931563
$ws.Range("C23").Value = 17700987132.42
$ws.Range("C24").Value = 0
$ws.Range("C25").Value = 2106581917.46
$ws.Range("C26").Value = 2366980990
$ws.Range("C27").Value = 52410709330.2726
$ws.Range("C28").Value = 112484268647.461
$ws.Range("C29").Value = $null
$ws.Range("C30").Value = $null
$ws.Range("C37").Value = 1528704399
$ws.Range("C38").Value = 1329297676
$ws.Range("C39").Value = 1914548429.67
$ws.Range("C40").Value = 21211340836.73
$ws.Range("C41").Value = 0
$ws.Range("C42").Value = 19739852323.73
$ws.Range("C43").Value = $null
$ws.Range("C45").Value = 908295196
$ws.Range("C48").Value = 2571353491
$ws.Range("C50").Value = 1183222441
$ws.Range("C51").Value = 9096048867.27738
$ws.Range("C52").Value = 30307389704.0074
$ws.Range("C53").Value = $null
$ws.Range("C54").Value = 28743629969
$ws.Range("C58").Value = 0
$ws.Range("C59").Value = -8280773107.3422
$ws.Range("C60").Value = 82583297861.65781
$ws.Range("C61").Value = -406418917.791066
$ws.Range("C62").Value = 82903158354.8102
$ws.Range("C63").Value = 112484268647.874
